# Rename the existing sheet, refresh its header row and append a sample
# data row, then add a second worksheet ("openAccountTest") with its own
# header row -- matches the "added excel utility with dataprovide" commit.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: AddMultipleCustomer -> addMultipleCustomerTest ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "addMultipleCustomerTest"

# Update header labels (row 1)
$ws1.Range("A1").Value = "First Name"
$ws1.Range("B1").Value = "Last Name"
$ws1.Range("C1").Value = "Post Code"

# Append a new sample data row (row 4)
$ws1.Range("A4").Value = "xyz"
$ws1.Range("B4").Value = "sample"
# Force text storage so the leading zeros in the post code survive
$ws1.Range("C4").NumberFormat = "@"
$ws1.Range("C4").Value = "00001"

# --- Sheet 2: new "openAccountTest" sheet placed after sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "openAccountTest"

$ws2.Range("A1").Value = "Customer"
$ws2.Range("B1").Value = "Currency"
$ws2.Range("A2").Select() | Out-Null

# Leave the first sheet active/selected, as in the source workbook
$ws1.Select() | Out-Null
$ws1.Range("D9").Select() | Out-Null
